# Reproduce the diagnostic.xlsx edit:
#  - B1 = 0, A2 = 0 : numeric cells with a bold font, thin box border,
#                     horizontally centered / top-aligned
#  - B2 = "disconnected_elements" : plain text (goes to sharedStrings.xml)
#  - sheet dimension grows from A1 to A1:B2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- values -----------------------------------------------------------
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# --- formatting ---------------------------------------------------------
# Build the full format (bold font + thin border on all sides + centered/top
# alignment) on B1 first. Doing every property mutation on a single Range
# keeps the style table minimal (one new cellXfs entry), then we clone that
# exact formatting onto A2 via copy/paste-special so no extra intermediate
# style combinations get created.
$r1 = $ws.Range("B1")
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108   # xlCenter
$r1.VerticalAlignment = -4160     # xlTop
$r1.Borders.LineStyle = 1         # xlContinuous
$r1.Borders.Weight = 2            # xlThin

$r1.Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

Write-Host "applied diagnostic formatting"
